{"js": "const replacements = [\n  [\"2024-09-01 Sunday\", \"2024-09-02 Monday\"],\n  [\"982\u00f78=122, 6\", \"165\u00f75=33, 0\"],\n  [\"367\u00f75=73, 2\", \"120\u00f73=40, 0\"],\n  [\"487\u00f78=60, 7\", \"529\u00f77=75, 4\"],\n  [\"426\u00f76=71, 0\", \"782\u00f76=130, 2\"],\n  [\"681\u00f73=227, 0\", \"538\u00f77=76, 6\"],\n  [\"893\u00f72=446, 1\", \"178\u00f72=89, 0\"],\n  [\"449\u00f79=49, 8\", \"185\u00f75=37, 0\"],\n  [\"469\u00f75=93, 4\", \"883\u00f76=147, 1\"],\n  [\"754\u00f76=125, 4\", \"131\u00f73=43, 2\"],\n  [\"164\u00f77=23, 3\", \"506\u00f72=253, 0\"],\n  [\"373\u00f75=74, 3\", \"117\u00f76=19, 3\"],\n  [\"923\u00f76=153, 5\", \"991\u00f72=495, 1\"],\n  [\"804\u00f79=89, 3\", \"936\u00f74=234, 0\"],\n  [\"152\u00f72=76, 0\", \"161\u00f74=40, 1\"],\n  [\"658\u00f72=329, 0\", \"728\u00f78=91, 0\"],\n  [\"390\u00f72=195, 0\", \"952\u00f74=238, 0\"],\n  [\"888\u00f77=126, 6\", \"787\u00f73=262, 1\"],\n  [\"949\u00f78=118, 5\", \"673\u00f77=96, 1\"],\n  [\"621\u00f76=103, 3\", \"435\u00f78=54, 3\"],\n  [\"988\u00f79=109, 7\", \"838\u00f77=119, 5\"],\n  [\"151\u00f73=50, 1\", \"150\u00f75=30, 0\"],\n  [\"854\u00f75=170, 4\", \"802\u00f74=200, 2\"],\n  [\"697\u00f75=139, 2\", \"741\u00f74=185, 1\"],\n  [\"783\u00f72=391, 1\", \"110\u00f75=22, 0\"],\n  [\"455\u00f78=56, 7\", \"194\u00f74=48, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  // Each snippet is unique in this document, so replace every hit found\n  // (normally exactly one) while keeping the run's own formatting intact.\n  for (const range of found.items) {\n    range.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-01 Sunday\", \"2024-09-02 Monday\"),\n    @(\"982\u00f78=122, 6\", \"165\u00f75=33, 0\"),\n    @(\"367\u00f75=73, 2\", \"120\u00f73=40, 0\"),\n    @(\"487\u00f78=60, 7\", \"529\u00f77=75, 4\"),\n    @(\"426\u00f76=71, 0\", \"782\u00f76=130, 2\"),\n    @(\"681\u00f73=227, 0\", \"538\u00f77=76, 6\"),\n    @(\"893\u00f72=446, 1\", \"178\u00f72=89, 0\"),\n    @(\"449\u00f79=49, 8\", \"185\u00f75=37, 0\"),\n    @(\"469\u00f75=93, 4\", \"883\u00f76=147, 1\"),\n    @(\"754\u00f76=125, 4\", \"131\u00f73=43, 2\"),\n    @(\"164\u00f77=23, 3\", \"506\u00f72=253, 0\"),\n    @(\"373\u00f75=74, 3\", \"117\u00f76=19, 3\"),\n    @(\"923\u00f76=153, 5\", \"991\u00f72=495, 1\"),\n    @(\"804\u00f79=89, 3\", \"936\u00f74=234, 0\"),\n    @(\"152\u00f72=76, 0\", \"161\u00f74=40, 1\"),\n    @(\"658\u00f72=329, 0\", \"728\u00f78=91, 0\"),\n    @(\"390\u00f72=195, 0\", \"952\u00f74=238, 0\"),\n    @(\"888\u00f77=126, 6\", \"787\u00f73=262, 1\"),\n    @(\"949\u00f78=118, 5\", \"673\u00f77=96, 1\"),\n    @(\"621\u00f76=103, 3\", \"435\u00f78=54, 3\"),\n    @(\"988\u00f79=109, 7\", \"838\u00f77=119, 5\"),\n    @(\"151\u00f73=50, 1\", \"150\u00f75=30, 0\"),\n    @(\"854\u00f75=170, 4\", \"802\u00f74=200, 2\"),\n    @(\"697\u00f75=139, 2\", \"741\u00f74=185, 1\"),\n    @(\"783\u00f72=391, 1\", \"110\u00f75=22, 0\"),\n    @(\"455\u00f78=56, 7\", \"194\u00f74=48, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nWrite-Output \"done\""}
